# Scheduled runner update: refresh Leve profit calculations (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 11 (Leve Item ID 5533)
$ws.Range("H11").Value = 691.4167
$ws.Range("I11").Value = 691.4167
$ws.Range("K11").Value = 691.4167
$ws.Range("M11").Value = -551.4167

# row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 800
$ws.Range("J32").Value = 800
$ws.Range("L32").Value = 800
$ws.Range("N32").Value = -1452

# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 3050
$ws.Range("J40").Value = 3050
$ws.Range("L40").Value = 3050
$ws.Range("N40").Value = -3400

# row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 687.75
$ws.Range("J55").Value = 1000.5
$ws.Range("L55").Value = 1000.5
$ws.Range("N55").Value = -1428.5

# row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 3928.4119
$ws.Range("I62").Value = 3599.5
$ws.Range("K62").Value = 3599.5
$ws.Range("M62").Value = -2975.5

# row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 3928.4119
$ws.Range("I65").Value = 3599.5
$ws.Range("K65").Value = 17997.5
$ws.Range("M65").Value = -14877.5

# row 95 (Leve Item ID 18200)
$ws.Range("H95").Value = 11000
$ws.Range("J95").Value = 11000
$ws.Range("L95").Value = 11000
$ws.Range("N95").Value = -16492

# row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 3447.55
$ws.Range("J129").Value = 4197.467
$ws.Range("L129").Value = 12592.401
$ws.Range("N129").Value = -22592.401

# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1250.3636
$ws.Range("I132").Value = 913.2353000000001
$ws.Range("K132").Value = 2739.7059
$ws.Range("M132").Value = -209.7058999999999

$ws = $wb.Worksheets.Item("ARM")
# row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 2182.2
$ws.Range("I74").Value = 2182.2
$ws.Range("K74").Value = 2182.2
$ws.Range("M74").Value = -1308.2

# row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 2182.2
$ws.Range("I77").Value = 2182.2
$ws.Range("K77").Value = 10911
$ws.Range("M77").Value = -6543

# row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2951.818
$ws.Range("I122").Value = 2558.875
$ws.Range("K122").Value = 7676.625
$ws.Range("M122").Value = -5226.625

# row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2043.7858
$ws.Range("I132").Value = 2051.1667
$ws.Range("K132").Value = 6153.500100000001
$ws.Range("M132").Value = -3623.500100000001

$ws = $wb.Worksheets.Item("BSM")
# row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 1626.6154
$ws.Range("I20").Value = 798.375
$ws.Range("K20").Value = 798.375
$ws.Range("M20").Value = -551.375

# row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 51
$ws.Range("I22").Value = 51
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 51
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 122
$ws.Range("N22").Value = ""

# row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 8600.083000000001
$ws.Range("I134").Value = 8577.888999999999
$ws.Range("K134").Value = 25733.667
$ws.Range("M134").Value = -23198.667

$ws = $wb.Worksheets.Item("CRP")
# row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 287.5
$ws.Range("I22").Value = 287.5
$ws.Range("K22").Value = 287.5
$ws.Range("M22").Value = 62.5

# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3906.6667
$ws.Range("I31").Value = 3532
$ws.Range("K31").Value = 3532
$ws.Range("M31").Value = -3237

# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3906.6667
$ws.Range("I34").Value = 3532
$ws.Range("K34").Value = 3532
$ws.Range("M34").Value = -3330

# row 47 (Leve Item ID 1920)
$ws.Range("H47").Value = 34999.5
$ws.Range("I47").Value = 19999
$ws.Range("K47").Value = 19999
$ws.Range("M47").Value = -19433

# row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1381
$ws.Range("I132").Value = 1381
$ws.Range("K132").Value = 4143
$ws.Range("M132").Value = -1613

# row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2761.8
$ws.Range("I134").Value = 1351.091
$ws.Range("K134").Value = 4053.273
$ws.Range("M134").Value = -1518.273

$ws = $wb.Worksheets.Item("CUL")
# row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 358.5
$ws.Range("I2").Value = 199
$ws.Range("J2").Value = 376.22223
$ws.Range("K2").Value = 1194
$ws.Range("L2").Value = 2257.33338
$ws.Range("M2").Value = -1081
$ws.Range("N2").Value = -2483.33338

# row 38 (Leve Item ID 4860)
$ws.Range("H38").Value = 32
$ws.Range("J38").Value = 15
$ws.Range("L38").Value = 45
$ws.Range("N38").Value = -739

# row 100 (Leve Item ID 19831)
$ws.Range("H100").Value = 3028
$ws.Range("J100").Value = 3028
$ws.Range("L100").Value = 9084
$ws.Range("N100").Value = -10706

# row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 1446
$ws.Range("J107").Value = 1630.625
$ws.Range("L107").Value = 4891.875
$ws.Range("N107").Value = -8731.875

# row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 2748.2856
$ws.Range("J131").Value = 2940
$ws.Range("L131").Value = 8820
$ws.Range("N131").Value = -18900

$ws = $wb.Worksheets.Item("GSM")
# row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 2999
$ws.Range("I113").Value = 2999
$ws.Range("K113").Value = 2999
$ws.Range("M113").Value = -829

# row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 8339393
$ws.Range("I122").Value = 13895766
$ws.Range("K122").Value = 41687298
$ws.Range("M122").Value = -41684848

# row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3391.4211
$ws.Range("I126").Value = 3209.4
$ws.Range("J126").Value = 3456.4285
$ws.Range("K126").Value = 9628.200000000001
$ws.Range("L126").Value = 10369.2855
$ws.Range("M126").Value = -7158.200000000001
$ws.Range("N126").Value = -15309.2855

$ws = $wb.Worksheets.Item("LTW")
# row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 3333.6924
$ws.Range("I40").Value = 2794.8333
$ws.Range("K40").Value = 2794.8333
$ws.Range("M40").Value = -2658.8333

# row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 3105.8667
$ws.Range("J46").Value = 4124.75
$ws.Range("L46").Value = 4124.75
$ws.Range("N46").Value = -4500.75

# row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 1634.5
$ws.Range("I93").Value = 1725.2858
$ws.Range("J93").Value = 999
$ws.Range("K93").Value = 1725.2858
$ws.Range("L93").Value = 999
$ws.Range("M93").Value = -477.2858000000001
$ws.Range("N93").Value = -3495

# row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 4145.7144
$ws.Range("I122").Value = 3503.3333
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 10509.9999
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -8059.999899999999
$ws.Range("N122").Value = -28900

# row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 5073.647
$ws.Range("I132").Value = 4883.533
$ws.Range("K132").Value = 14650.599
$ws.Range("M132").Value = -12120.599

# row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 2657.625
$ws.Range("I136").Value = 2534.4666
$ws.Range("K136").Value = 7603.399800000001
$ws.Range("M136").Value = -5053.399800000001

$ws = $wb.Worksheets.Item("WVR")
# row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 1257.4
$ws.Range("I122").Value = 1257.4
$ws.Range("K122").Value = 3772.2
$ws.Range("M122").Value = -1322.2

# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1054.2858
$ws.Range("I132").Value = 1055
$ws.Range("J132").Value = 1050
$ws.Range("K132").Value = 3165
$ws.Range("L132").Value = 3150
$ws.Range("M132").Value = -635
$ws.Range("N132").Value = -8210

